$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.331.04'
$ws.Range("E2").Value = '  +4.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.580.16'
$ws.Range("E3").Value = '  +0.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -1.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.54'
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("E7").Value = '  -1.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.63'
$ws.Range("E8").Value = '  +7.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.252'
$ws.Range("E9").Value = '  +0.92%  '

$ws.Range("E10").Value = '  -0.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("E11").Value = '  +1.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.805.81'
$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.570.02'
$ws.Range("E13").Value = '  -0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.76'
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.305.55'
$ws.Range("E16").Value = '  +4.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.82'
$ws.Range("E17").Value = '  +2.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.38'
$ws.Range("E18").Value = '  +7.29%  '

$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("E22").Value = '  -0.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.36'
$ws.Range("E23").Value = '  +1.16%  '

$ws.Range("E24").Value = '  -0.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.55'
$ws.Range("E25").Value = '  -1.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.28'
$ws.Range("E26").Value = '  +1.04%  '

$ws.Range("E27").Value = '  -1.18%  '

$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.14'
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.14'
$ws.Range("E33").Value = '  -1.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.419.05'
$ws.Range("E34").Value = '  -2.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").Value = '  -1.15%  '

$ws.Range("E36").Value = '  -5.15%  '

$ws.Range("E37").Value = '  -1.37%  '

$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.52'
$ws.Range("E39").Value = '  +7.22%  '

$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.67'
$ws.Range("E43").Value = '  -2.46%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.83'
$ws.Range("E44").Value = '  +5.24%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.973'
$ws.Range("E45").Value = '  -2.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.20'
$ws.Range("E46").Value = '  -0.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.716.96'
$ws.Range("E47").Value = '  +0.72%  '

$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("E50").Value = '  +0.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.37'
$ws.Range("E51").Value = '  +16.11%  '

